# Insert a new weekly price row at row 40, shifting all existing data rows
# (previously rows 40-128) down by one (to rows 41-129).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(40).Insert()

# Populate the newly inserted row 40 with the new data point.
$ws.Cells.Item(40, 1).Value = 7
$ws.Cells.Item(40, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(40, 3).Value = "Ñuble"
$ws.Cells.Item(40, 4).Value = 45238
$ws.Cells.Item(40, 5).Value = 16
$ws.Cells.Item(40, 6).Value = 100112001
$ws.Cells.Item(40, 7).Value = "Berenjena"
$ws.Cells.Item(40, 8).Value = "Sin especificar"
$ws.Cells.Item(40, 9).Value = "Primera"
$ws.Cells.Item(40, 10).Value = 60
$ws.Cells.Item(40, 11).Value = 12000
$ws.Cells.Item(40, 12).Value = 12000
$ws.Cells.Item(40, 13).Value = 12000
$ws.Cells.Item(40, 14).Value = "$/caja 60 unidades"
$ws.Cells.Item(40, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(40, 16).Value = 200
$ws.Cells.Item(40, 17).Value = 60
$ws.Cells.Item(40, 18).Value = "Hortaliza"
